$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates: (row, col, value, forceText)
$data = @(
    @(2, "D", "29.912.26", $false),
    @(2, "E", "  +0.07%  ", $false),
    @(3, "D", "1.874.84", $false),
    @(3, "E", "  -0.89%  ", $false),
    @(4, "E", "  -0.06%  ", $false),
    @(5, "D", "0.7414", $true),
    @(5, "E", "  -4.02%  ", $false),
    @(6, "D", "242.46", $true),
    @(6, "E", "  -0.47%  ", $false),
    @(7, "E", "  -0.05%  ", $false),
    @(8, "D", "0.3150", $true),
    @(8, "E", "  +0.82%  ", $false),
    @(9, "D", "0.07160", $true),
    @(9, "E", "  -1.03%  ", $false),
    @(10, "D", "24.71", $true),
    @(10, "E", "  -3.67%  ", $false),
    @(11, "D", "0.08406", $true),
    @(11, "E", "  -3.59%  ", $false),
    @(12, "D", "0.7514", $true),
    @(12, "E", "  -2.49%  ", $false),
    @(13, "D", "5.420", $true),
    @(13, "E", "  +0.38%  ", $false),
    @(14, "D", "1.839.89", $false),
    @(14, "E", "  -13.40%  ", $false),
    @(15, "D", "92.58", $true),
    @(15, "E", "  -1.81%  ", $false),
    @(16, "B", "Uniswap", $false),
    @(16, "C", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", $false),
    @(16, "D", "6.110", $true),
    @(16, "E", "  -1.63%  ", $false),
    @(17, "B", "WrappedBTC", $false),
    @(17, "C", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", $false),
    @(17, "D", "29.909.08", $false),
    @(17, "E", "  +0.05%  ", $false),
    @(18, "D", "13.59", $true),
    @(18, "E", "  -2.30%  ", $false),
    @(19, "D", "243.42", $true),
    @(19, "E", "  -0.73%  ", $false),
    @(20, "D", "0.000007816", $true),
    @(20, "E", "  -0.70%  ", $false),
    @(21, "D", "0.9987", $true),
    @(21, "E", "  -0.16%  ", $false),
    @(22, "D", "2.123.23", $false),
    @(22, "E", "  -12.35%  ", $false),
    @(23, "D", "7.986", $true),
    @(23, "E", "  -2.30%  ", $false),
    @(24, "D", "0.9996", $true),
    @(24, "E", "  -0.15%  ", $false),
    @(25, "D", "0.1551", $true),
    @(25, "E", "  -2.49%  ", $false),
    @(26, "D", "9.297", $true),
    @(26, "E", "  -2.28%  ", $false),
    @(27, "D", "165.33", $true),
    @(27, "E", "  +1.80%  ", $false),
    @(28, "E", "  -0.98%  ", $false),
    @(29, "E", "  -0.40%  ", $false),
    @(30, "D", "1.485", $true),
    @(30, "E", "  +3.81%  ", $false),
    @(31, "D", "4.614", $true),
    @(31, "E", "  +2.21%  ", $false),
    @(32, "D", "1.527", $true),
    @(32, "E", "  -1.14%  ", $false),
    @(33, "D", "4.263", $true),
    @(33, "E", "  +3.49%  ", $false),
    @(34, "D", "0.05331", $true),
    @(34, "E", "  -1.99%  ", $false),
    @(35, "D", "1.238", $true),
    @(35, "E", "  -0.74%  ", $false),
    @(36, "D", "0.7546", $true),
    @(36, "E", "  +0.47%  ", $false),
    @(37, "D", "0.9966", $true),
    @(37, "E", "  -0.90%  ", $false),
    @(38, "D", "2.697", $true),
    @(38, "E", "  +0.07%  ", $false),
    @(39, "D", "0.01950", $true),
    @(39, "E", "  -1.65%  ", $false),
    @(40, "D", "2.752", $true),
    @(40, "E", "  -1.19%  ", $false),
    @(41, "D", "0.4494", $true),
    @(41, "E", "  -0.44%  ", $false),
    @(42, "D", "1.110.38", $false),
    @(42, "E", "  +1.02%  ", $false),
    @(43, "D", "6.058", $true),
    @(43, "E", "  -0.67%  ", $false),
    @(44, "D", "72.19", $true),
    @(44, "E", "  -1.63%  ", $false),
    @(45, "D", "0.8567", $true),
    @(45, "E", "  +0.18%  ", $false),
    @(46, "E", "  +0.13%  ", $false),
    @(47, "D", "103.13", $true),
    @(47, "E", "  -0.36%  ", $false),
    @(48, "D", "7.663", $true),
    @(48, "E", "  +0.66%  ", $false),
    @(49, "D", "3.085", $true),
    @(49, "E", "  +2.98%  ", $false),
    @(50, "D", "1.839", $true),
    @(50, "E", "  -2.42%  ", $false),
    @(51, "D", "2.020.78", $false),
    @(51, "E", "  -7.53%  ", $false)
)

foreach ($item in $data) {
    $r = $item[0]
    $col = $item[1]
    $val = $item[2]
    $forceText = $item[3]
    $addr = "$col$r"
    $rng = $ws.Range($addr)
    if ($forceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $val
}

Write-Host "Applied $($data.Count) cell updates"